$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 header labels: drop the "(n)" footnote runs, keep plain text.
# Column D header changes to a brand-new label ("Hoa hong thuc hien theo Y/cau"
# is replaced by "Hoa hong NV ho tro").
$ws.Range("A5").Value = "Mã nhân viên"
$ws.Range("B5").Value = "Tên nhân viên"
$ws.Range("C5").Value = "Hoa hồng thực hiện"
$ws.Range("D5").Value = "Hoa hồng NV hỗ trợ"
$ws.Range("E5").Value = "Hoa hồng tư vấn"
$ws.Range("F5").Value = "Hoa hồng bán gói dịch vụ"
$ws.Range("G5").Value = "Tổng"

# Totals-row label: drop the trailing colon.
$ws.Range("A30").Value = "Tổng cộng"

# The merged totals cell (A30:B30) used to carry a left border on A30 only;
# simplify it so the merged cell just keeps its top/bottom rule (matching B30).
$ws.Range("A30").Borders.Item(7).LineStyle = -4142

# Selection moved to A9 when the file was last saved.
$ws.Range("A9").Select()
